# Update generated output values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 8457
$wsExpo.Range("F3").Value = 8113
$wsExpo.Range("F10").Value = 194
$wsExpo.Range("F11").Value = 244
$wsExpo.Range("F12").Value = 733
$wsExpo.Range("F13").Value = 189
$wsExpo.Range("F14").Value = 4079
$wsExpo.Range("F16").Value = 72
$wsExpo.Range("F19").Value = 146
$wsExpo.Range("F20").Value = 105

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 8457
$wsAll.Range("F3").Value = 8113
$wsAll.Range("F10").Value = 194
$wsAll.Range("F11").Value = 244
$wsAll.Range("F12").Value = 733
$wsAll.Range("F13").Value = 189
$wsAll.Range("F14").Value = 4080
$wsAll.Range("F16").Value = 72
$wsAll.Range("F19").Value = 146
$wsAll.Range("F20").Value = 105
